$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The blank spacer row that used to sit between row 2 and the "Handle DAQmx
# errors..." entry is removed; everything below it shifts up by one row.
$ws.Rows("3").Delete()

# Fill in the newly freed-up last row of the table (row 7) with the new
# "task controller" issue described in the commit message.
$ws.Range("A7").Value = "7/20/2015"
$ws.Range("B7").Value = "high"
$ws.Range("C7").Value = "task controller"
$ws.Range("D7").Value = "if task controller is set to iterate in parallel with child TCs then it receives an iteration event while it is in a done state. This error was seen before but it was corrected. It seems that it was corrected only for the case when the task controller is set to iterate before child TCs. In this case everything works fine. To reproduce the error, set up an UITC, add the scan engine as child and set it to iterate in parallel with child TCs. Add as children to the scan engine the DAQdev and pockells modules. Then set the UITC to do multiple iterations."
$ws.Rows("7").RowHeight = 90

# Update the view's selection to the new range used while entering this data.
$ws.Range("D8:D13").Select()
